$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column S with the 2022 data by copying column R's formatting
# (same pattern Excel uses when a new year's column is added next to the
# previous one) and then filling in the new values.
$ws.Range("R1:R35").Copy()
$ws.Range("S1:S35").Insert(-4161)

# New 2022 figures (column S), row by row.
$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36

$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17

$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6

$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2

$ws.Range("S17").Value = "-"
$ws.Range("S18").Value = "-"
$ws.Range("S19").Value = "-"

$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5

$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"

$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6

$ws.Range("S29").Value = "-"
$ws.Range("S30").Value = "-"
$ws.Range("S31").Value = "-"

$ws.Range("S32").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("S34").Value = "-"

# Update the active selection to reflect where the user ended up.
$ws.Range("S3").Select()
